# Fruta / hortaliza, semanal
# A new weekly observation was inserted at row 579 ("Femacal de La Calera" /
# "Poroto verde" price list). Inserting the row pushes every following
# record (old rows 579-656) down by one (new rows 580-657), which is why
# the sheet's used range grows from A1:R656 to A1:R657.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 579; everything below (old 579..656)
# shifts down to 580..657 automatically.
$ws.Rows.Item(579).Insert()

# The new row reuses the same "shape" as its neighbour (now at row 580,
# which used to be row 579) for every column that does not change, so
# copy that row across first and then patch just the columns called out
# by the diff.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(579, $col).Value = $ws.Cells.Item(580, $col).Value()
}

# Patch the cells that actually differ for the new record.
$ws.Cells.Item(579, 4).Value  = 45154              # D579 Fecha
$ws.Cells.Item(579, 8).Value  = "Sin especificar"  # H579 Variedad
$ws.Cells.Item(579, 11).Value = 34000              # K579 Precio minimo
$ws.Cells.Item(579, 12).Value = 35000              # L579 Precio maximo
$ws.Cells.Item(579, 13).Value = 34521              # M579 Precio promedio ponderado
$ws.Cells.Item(579, 15).Value = "Perú"             # O579 Origen
$ws.Cells.Item(579, 16).Value = 1381               # P579 Precio $/Kg
